$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Team Month")

# Update the "Team" column labels (B2:B5) to their new values.
$ws.Range("B2").Value = "Xinghao_s2l"
$ws.Range("B3").Value = "EDN_S2l"
$ws.Range("B4").Value = "Cpu_s2l"
$ws.Range("B5").Value = "Kurni_s2l"

# Update Weekly Pending Total (C) and Repayment (D) figures for every row.
# Column E (Recovery rate) already holds D/C formulas and recalculates automatically.
$ws.Range("C2").Value = 2219432362
$ws.Range("D2").Value = 231288930

$ws.Range("C3").Value = 1487395776
$ws.Range("D3").Value = 150411420

$ws.Range("C4").Value = 1503074041
$ws.Range("D4").Value = 149033593

$ws.Range("C5").Value = 5424927775
$ws.Range("D5").Value = 530588991

$ws.Range("C6").Value = 5500251169
$ws.Range("D6").Value = 482558848

$ws.Range("C7").Value = 2889911911
$ws.Range("D7").Value = 232381599

$ws.Range("C8").Value = 5383097040
$ws.Range("D8").Value = 430019048

# Move the active selection to A2 (matches the saved view state in the file).
$ws.Range("A2").Select() | Out-Null
